$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1) - rows 2..23
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 1253
$ws1.Range("F3").Value = 17106
$ws1.Range("F4").Value = 52
$ws1.Range("F5").Value = 1667
$ws1.Range("G7").Value = 49.9
$ws1.Range("F8").Value = 1045
$ws1.Range("F9").Value = 402
$ws1.Range("F11").Value = 137
$ws1.Range("F12").Value = 11891
$ws1.Range("F15").Value = 11583
$ws1.Range("F18").Value = 59
$ws1.Range("F19").Value = 414
$ws1.Range("F21").Value = 921
$ws1.Range("F22").Value = 347
$ws1.Range("F23").Value = 155

# Sheet "全部类型" (sheet4) - rows shifted by +2 after row 11
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 1253
$ws4.Range("F3").Value = 17106
$ws4.Range("F4").Value = 52
$ws4.Range("F5").Value = 1667
$ws4.Range("G7").Value = 49.9
$ws4.Range("F8").Value = 1045
$ws4.Range("F9").Value = 402
$ws4.Range("F11").Value = 137
$ws4.Range("F14").Value = 11891
$ws4.Range("F17").Value = 11583
$ws4.Range("F20").Value = 59
$ws4.Range("F21").Value = 414
$ws4.Range("F23").Value = 921
$ws4.Range("F24").Value = 347
$ws4.Range("F25").Value = 155
